# Auto-generated edit script applying the cryptos price-list diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "62.932.73"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3
$ws.Range("D3").Value = "2.463.70"
$ws.Range("E3").Value = "  +0.76%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.70"
$ws.Range("E5").Value = "  -0.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.57"
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").Value = "2.463.72"
$ws.Range("E9").Value = "  +0.82%  "

# Row 10
$ws.Range("E10").Value = "  +1.31%  "

# Row 11
$ws.Range("E11").Value = "  +1.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.27"
$ws.Range("E12").Value = "  +0.79%  "

# Row 13
$ws.Range("E13").Value = "  +0.98%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.05"
$ws.Range("E14").Value = "  +2.14%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("E15").Value = "  +0.18%  "

# Row 16
$ws.Range("D16").Value = "2.911.08"
$ws.Range("E16").Value = "  +0.77%  "

# Row 17
$ws.Range("D17").Value = "62.911.74"
$ws.Range("E17").Value = "  +0.56%  "

# Row 18
$ws.Range("D18").Value = "2.469.73"
$ws.Range("E18").Value = "  +1.73%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.97"
$ws.Range("E19").Value = "  +2.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.04"
$ws.Range("E20").Value = "  +1.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.55"
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("E22").Value = "  +11.54%  "

# Row 24
$ws.Range("B24").Value = "Aptos"
$ws.Range("C24").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.39"
$ws.Range("E24").Value = "  +22.62%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.87"
$ws.Range("E26").Value = "  +0.54%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "652.06"
$ws.Range("E27").Value = "  +1.14%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0985"
$ws.Range("E28").Value = "  +0.73%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -14.07%  "

# Row 31
$ws.Range("E31").Value = "  +3.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("E32").Value = "  -1.83%  "

# Row 33
$ws.Range("E33").Value = "  -0.56%  "

# Row 34
$ws.Range("E34").Value = "  -3.90%  "

# Row 36
$ws.Range("E36").Value = "  +4.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("E37").Value = "  +0.67%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.370"
$ws.Range("E38").Value = "  -0.69%  "

# Row 39
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.36"
$ws.Range("E39").Value = "  -0.54%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.43"
$ws.Range("E40").Value = "  -0.68%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.74"
$ws.Range("E41").Value = "  +0.82%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.78"
$ws.Range("E42").Value = "  +2.89%  "

# Row 43
$ws.Range("E43").Value = "  -1.17%  "

# Row 44
$ws.Range("D44").Value = "0.0₆0318"
$ws.Range("E44").Value = "  -63.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.88"
$ws.Range("E46").Value = "  +6.79%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.25"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.59"
$ws.Range("E48").Value = "  +0.13%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.49"
$ws.Range("E49").Value = "  -0.52%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.608"
$ws.Range("E50").Value = "  +0.53%  "

# Row 51
$ws.Range("E51").Value = "  +0.22%  "
